$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bab 4")

$ws.Range("H1").Value = "Tabel 4.2.5"

# P1 and W1 keep the "Tabel" prefix underlined (inherited from the cell
# style) while only the number suffix is rendered without the underline -
# mirrors the rich-text run split already present in the source file.
$ws.Range("P1").Value = "Tabel 4.2.6."
$ws.Range("P1").Characters(6, 7).Font.Underline = $false

$ws.Range("W1").Value = "Tabel 4.2.7."
$ws.Range("W1").Characters(6, 7).Font.Underline = $false

$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan {kec}. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan {kec}, 2021"
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan {kec}, 2021"
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan {kec}, 2021"

$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in {kec} Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village {kec} Subdistrict, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in {kec} Subdistrict, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in {kec} Subdistrict, 2021"
